$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "DONE"
$ws.Range("E2").Value = "DONE"

$ws.Range("D3").Select()
